# Rename the field/value mapping sheets for excel/csv consistency,
# and make the "values" sheet the active/selected sheet.

$wb = $excel.ActiveWorkbook

$fieldsSheet = $wb.Worksheets.Item("field_mapping")
$fieldsSheet.Name = "fields"

$valuesSheet = $wb.Worksheets.Item("value_mapping")
$valuesSheet.Name = "values"

$valuesSheet.Select()
